$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @(0.3464964993005633, 9.226618575922256, 2938.103010863317, 6.48142807727062, 2954.157554015811)
    3 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.034748368925986)
    4 = @(0.06328177979961902, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 0.768386970581898)
    5 = @(1.505614041169197, 10990084.13351303, 16.98373111632243, 5548678842208.939, 5548689832311.562)
    6 = @(0.001754667048134761, 0.3375848360084654, 3.082599426703578, 6.48142807727062, 9.903367007030798)
    7 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    8 = @(0.006876353814593728, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 1.271808613593681)
    9 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
